$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New problem rows (156-185) for the "Сравнение по модулю" subtopic under "ТЧ",
# added per commit "added competetive_points and work_points".
$newRows = @(
@("ТЧ","Сравнение по модулю",0.0333,0.0667,2,1,1,5,2,3,3,2,1,1,0,1),
@("ТЧ","Сравнение по модулю",0.0667,0.0667,2,1,1,5,2,3,3,2,1,1,0,1),
@("ТЧ","Сравнение по модулю",0.1,0.0667,2,1,1,5,2,3,3,2,1,1,0,1),
@("ТЧ","Сравнение по модулю",0.1333,0.1333,2,1,1,5,2,3,3,2,1,3,0,1),
@("ТЧ","Сравнение по модулю",0.1667,0.1333,0,1,3,0,2,0,3,4,1,5,0,1),
@("ТЧ","Сравнение по модулю",0.2,0.2,2,2,3,5,4,5,3,4,1,3,0,5),
@("ТЧ","Сравнение по модулю",0.2333,0.2,0,2,3,5,4,0,5,4,3,3,0,5),
@("ТЧ","Сравнение по модулю",0.2667,0.2,0,0,3,5,5,0,0,4,3,0,0,5),
@("ТЧ","Сравнение по модулю",0.3,0.2667,2,2,3,5,4,5,5,4,3,3,0,5),
@("ТЧ","Сравнение по модулю",0.3333,0.2667,0,0,3,0,0,5,5,0,2,5,0,5),
@("ТЧ","Сравнение по модулю",0.3667,0.3333,2,4,3,5,0,5,5,2,3,4,0,5),
@("ТЧ","Сравнение по модулю",0.4,0.3333,4,4,3,5,0,5,5,5,3,4,0,5),
@("ТЧ","Сравнение по модулю",0.4333,0.3333,5,0,3,5,0,0,5,5,3,0,0,5),
@("ТЧ","Сравнение по модулю",0.4667,0.4,2,4,3,5,0,5,5,2,3,0,0,5),
@("ТЧ","Сравнение по модулю",0.5,0.4,4,4,3,5,0,5,5,5,2,5,0,5),
@("ТЧ","Сравнение по модулю",0.5333,0.4,5,4,3,0,0,0,0,0,2,0,0,5),
@("ТЧ","Сравнение по модулю",0.5667,0.4667,5,0,2,0,0,0,0,0,0,0,0,5),
@("ТЧ","Сравнение по модулю",0.6,0.5333,0,0,4,0,0,0,0,0,0,0,0,5),
@("ТЧ","Сравнение по модулю",0.6333,0.6,4,4,2,0,5,0,0,5,3,0,0,5),
@("ТЧ","Сравнение по модулю",0.6667,0.6,5,4,2,0,5,0,0,5,3,0,0,5),
@("ТЧ","Сравнение по модулю",0.7,0.6,5,5,5,0,5,0,0,0,3,0,0,5),
@("ТЧ","Сравнение по модулю",0.7333,0.6667,0,0,5,0,0,0,0,0,0,0,0,5),
@("ТЧ","Сравнение по модулю",0.7667,0.6667,0,0,5,0,0,0,0,0,0,0,0,5),
@("ТЧ","Сравнение по модулю",0.8,0.7333,0,0,5,0,0,0,0,0,0,0,0,0),
@("ТЧ","Сравнение по модулю",0.8333,0.7333,0,0,5,0,0,0,0,0,0,0,0,5),
@("ТЧ","Сравнение по модулю",0.8667,0.7333,0,0,5,0,0,0,0,0,0,0,0,5),
@("ТЧ","Сравнение по модулю",0.9,0.8,0,0,5,0,0,0,0,0,0,0,0,5),
@("ТЧ","Сравнение по модулю",0.9333,0.8667,0,0,5,0,0,0,0,0,0,0,0,5),
@("ТЧ","Сравнение по модулю",0.9667,0.9333,0,0,5,0,0,0,0,0,5,0,0,5),
@("ТЧ","Сравнение по модулю",1,1,0,0,0,0,0,0,0,0,0,0,0,5)
)

$startRow = 156
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $vals = $newRows[$i]
    for ($c = 0; $c -lt $vals.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $vals[$c]
    }
}

# Match the saved view: scrolled down and a range selected below the new data.
$ws.Range("Q156:T185").Select()
